$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ArrayInfo (C2): collapse the old multi-line 6x6 block array into a single
# flat, comma-separated row-major array (fixes MonsterPosList's
# ParseVector2int parsing bug referenced in the commit message).
$ws.Range("C2").Value = "{0,0,0,0,0,0,0,1,0,0,0,0,0,1,0,0,0,0,0,1,1,1,1,0,0,0,0,0,1,0,0,0,0,0,0,0}"

# MonsterNameList (F2): was empty, now spawns a BlueMonster.
$ws.Range("F2").Value = "{BlueMonster}"

# The row shrinks now that ArrayInfo is a single wrapped line instead of six.
$ws.Rows.Item(2).RowHeight = 66

# Move the active selection.
[void]$ws.Range("E3").Select()
